$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 7 (Ano 2025) figures as per latest data refresh
$ws.Range("B7").Value = 3705484.21
$ws.Range("C7").Value = -16.60096649392246
$ws.Range("D7").Value = 3275
$ws.Range("E7").Value = 3275
$ws.Range("F7").Value = 1131.445560305343
$ws.Range("G7").Value = 20.60391532359791
